$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Write brand-new text values first, in the exact order they must be
# --- appended to the shared-strings table (bottom-up, matching how the
# --- original author inserted rows working from the end of the list
# --- upward): "Dockspace Menu" (row 22) -> "ShaderProgram" (row 20) ->
# --- "ImGui Rendering" (row 14).
$ws.Range("A22").Value = "Dockspace Menu"
$ws.Range("A20").Value = "ShaderProgram"
$ws.Range("A14").Value = "ImGui Rendering"

# --- Now lay out the rest of column A (task names) for rows 13-22 ---
$ws.Range("A13").Value = "Docking"
$ws.Range("A15").Value = "Buffer"
$ws.Range("A16").Value = "Vertex Buffer"
$ws.Range("A17").Value = "Index Buffer"
$ws.Range("A18").Value = "Vertex Array Object"
$ws.Range("A19").Value = "Shaders"
$ws.Range("A21").Value = "Timer"

# --- Column B status values for rows 13-22 ---
$ws.Range("B13").Value = "Yes"
$ws.Range("B14").Value = "No"
$ws.Range("B15").Value = "No"
$ws.Range("B16").Value = "No"
$ws.Range("B17").Value = "No"
$ws.Range("B18").Value = "No"
$ws.Range("B19").Value = "Yes"
$ws.Range("B20").Value = "Yes"
$ws.Range("B21").Value = "No"
$ws.Range("B22").Value = "No"

# --- Update the view: scroll so row 10 is at the top and select B14 ---
$ws.Activate()
try {
    $excel.ActiveWindow.ScrollRow = 10
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
$ws.Range("B14").Select()
